$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.720.00"
$ws.Range("E2").Value = "'  +0.38%  "
$ws.Range("D3").Value = "'1.600.72"
$ws.Range("E3").Value = "'  +0.25%  "
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("D5").Value = "'211.31"
$ws.Range("E5").Value = "'  +0.13%  "
$ws.Range("E6").Value = "'  -0.49%  "
$ws.Range("E8").Value = "'  +0.27%  "
$ws.Range("E9").Value = "'  +0.33%  "
$ws.Range("E10").Value = "'  +1.09%  "
$ws.Range("E11").Value = "'  +0.83%  "
$ws.Range("D12").Value = "'1.825.58"
$ws.Range("E12").Value = "'  +0.25%  "
$ws.Range("D13").Value = "'1.611.65"
$ws.Range("E13").Value = "'  +1.66%  "
$ws.Range("E14").Value = "'  +0.49%  "
$ws.Range("E15").Value = "'  +0.23%  "
$ws.Range("D16").Value = "'65.19"
$ws.Range("E16").Value = "'  +0.29%  "
$ws.Range("D17").Value = "'26.696.91"
$ws.Range("E18").Value = "'  +0.83%  "
$ws.Range("D19").Value = "'210.75"
$ws.Range("E19").Value = "'  +1.13%  "
$ws.Range("E20").Value = "'  +2.54%  "
$ws.Range("D21").Value = "'1.01"
$ws.Range("E21").Value = "'  +0.18%  "
$ws.Range("E22").Value = "'  +0.86%  "
$ws.Range("E23").Value = "'  +0.08%  "
$ws.Range("E24").Value = "'  +0.91%  "
$ws.Range("D25").Value = "'143.79"
$ws.Range("E25").Value = "'  -0.97%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "'  +0.19%  "
$ws.Range("E27").Value = "'  -0.26%  "
$ws.Range("E28").Value = "'  -0.82%  "
$ws.Range("D29").Value = "'15.38"
$ws.Range("E29").Value = "'  +0.87%  "
$ws.Range("E30").Value = "'  +1.22%  "
$ws.Range("E31").Value = "'  -0.28%  "
$ws.Range("D33").Value = "'2.98"
$ws.Range("E33").Value = "'  +1.62%  "
$ws.Range("D34").Value = "'1.296.17"
$ws.Range("E35").Value = "'  +0.75%  "
$ws.Range("E36").Value = "'  -1.39%  "
$ws.Range("E37").Value = "'  +1.10%  "
$ws.Range("E38").Value = "'  +20.08%  "
$ws.Range("E39").Value = "'  -0.49%  "
$ws.Range("E40").Value = "'  -1.71%  "
$ws.Range("E41").Value = "'  -1.12%  "
$ws.Range("E42").Value = "'  -0.13%  "
$ws.Range("D43").Value = "'0.783"
$ws.Range("E43").Value = "'  -0.27%  "
$ws.Range("D44").Value = "'63.26"
$ws.Range("E44").Value = "'  -1.16%  "
$ws.Range("D45").Value = "'1.737.66"
$ws.Range("E45").Value = "'  +0.26%  "
$ws.Range("D46").Value = "'91.15"
$ws.Range("E46").Value = "'  +1.21%  "
$ws.Range("E47").Value = "'  -2.60%  "
$ws.Range("B48").Value = "'BabyDogeCoin"
$ws.Range("C48").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.0₆0105"
$ws.Range("E48").Value = "'  -1.05%  "
$ws.Range("B49").Value = "'Algorand"
$ws.Range("C49").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.102"
$ws.Range("E49").Value = "'  -0.34%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0518"
$ws.Range("E50").Value = "'  +1.99%  "
$ws.Range("B51").Value = "'USDD"
$ws.Range("C51").Value = "'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "'  +0.08%  "
